# Update cryptocurrency price/volume data and reorder a block of coin rows
# as published by the "Updated symbol list" GitHub Actions workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell below is stored in the workbook as plain text (inline/shared string),
# even though many of the values look numeric (e.g. "314.72") or percentage-like
# (e.g. "1.33%"). Assigning a bare numeric-looking string via .Value would cause
# Excel to auto-convert it into a real number, which would change both the cell
# type and its visual formatting. To avoid that, for every target cell we:
#   1. Force the NumberFormat to Text ("@") so Excel stores the literal text,
#   2. Assign the value,
#   3. Reset the cell style back to "Normal" so no stray number-format style
#      is left behind (matching the original workbook's unstyled text cells).

$updates = @(
    @{ Cell = 'D2'; Value = '314.72' },
    @{ Cell = 'E2'; Value = '1.33%' },
    @{ Cell = 'D3'; Value = '40.98' },
    @{ Cell = 'E3'; Value = '-0.80%' },
    @{ Cell = 'D4'; Value = '5.122' },
    @{ Cell = 'E4'; Value = '-0.08%' },
    @{ Cell = 'D5'; Value = '0.07633' },
    @{ Cell = 'E5'; Value = '-0.84%' },
    @{ Cell = 'B6'; Value = 'FTXToken' },
    @{ Cell = 'C6'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' },
    @{ Cell = 'D6'; Value = '1.687' },
    @{ Cell = 'E6'; Value = '3.70%' },
    @{ Cell = 'B7'; Value = 'MXToken' },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D7'; Value = '0.9357' },
    @{ Cell = 'E7'; Value = '1.30%' },
    @{ Cell = 'B8'; Value = 'BTSEToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' },
    @{ Cell = 'D8'; Value = '2.425' },
    @{ Cell = 'E8'; Value = '-1.74%' },
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D9'; Value = '0.1249' },
    @{ Cell = 'E9'; Value = '3.18%' },
    @{ Cell = 'B10'; Value = 'WazirX' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'D10'; Value = '0.1828' },
    @{ Cell = 'E10'; Value = '-0.05%' },
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D11'; Value = '0.09026' },
    @{ Cell = 'E11'; Value = '-1.85%' },
    @{ Cell = 'B12'; Value = 'BitrueCoin' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'D12'; Value = '0.04138' },
    @{ Cell = 'E12'; Value = '-1.84%' },
    @{ Cell = 'B13'; Value = 'BitMartToken' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'D13'; Value = '0.1055' },
    @{ Cell = 'E13'; Value = '0.43%' },
    @{ Cell = 'B14'; Value = 'BitForexToken' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'D14'; Value = '0.001264' },
    @{ Cell = 'E14'; Value = '0.86%' },
    @{ Cell = 'B15'; Value = 'TigerCash' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'D15'; Value = '0.005856' },
    @{ Cell = 'E15'; Value = '2.58%' },
    @{ Cell = 'B16'; Value = 'UpBots' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt' },
    @{ Cell = 'D16'; Value = '0.007491' },
    @{ Cell = 'E16'; Value = '1,897.31%' },
    @{ Cell = 'B17'; Value = 'LEO' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'D17'; Value = '3.354' },
    @{ Cell = 'E17'; Value = '0.08%' },
    @{ Cell = 'B18'; Value = 'GateToken' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'D18'; Value = '4.334' },
    @{ Cell = 'E18'; Value = '0.60%' },
    @{ Cell = 'D19'; Value = '0.3359' },
    @{ Cell = 'E19'; Value = '1.68%' },
    @{ Cell = 'D20'; Value = '8.432' },
    @{ Cell = 'E20'; Value = '21.57%' },
    @{ Cell = 'E21'; Value = '-2.97%' },
    @{ Cell = 'D23'; Value = '0.04040' },
    @{ Cell = 'E23'; Value = '-0.20%' },
    @{ Cell = 'D24'; Value = '0.001266' },
    @{ Cell = 'E24'; Value = '0.43%' },
    @{ Cell = 'D25'; Value = '0.004046' },
    @{ Cell = 'E25'; Value = '-1.41%' },
    @{ Cell = 'D26'; Value = '0.0001275' },
    @{ Cell = 'E26'; Value = '0.45%' },
    @{ Cell = 'D38'; Value = '0.02479' },
    @{ Cell = 'E38'; Value = '0.36%' },
    @{ Cell = 'D39'; Value = '0.05198' },
    @{ Cell = 'E39'; Value = '-1.33%' },
    @{ Cell = 'D40'; Value = '0.007787' },
    @{ Cell = 'E40'; Value = '-0.52%' },
    @{ Cell = 'D41'; Value = '0.1297' },
    @{ Cell = 'E41'; Value = '-1.39%' },
    @{ Cell = 'D42'; Value = '0.007370' },
    @{ Cell = 'E42'; Value = '8.58%' },
    @{ Cell = 'D43'; Value = '0.002169' },
    @{ Cell = 'E43'; Value = '17.68%' },
    @{ Cell = 'D44'; Value = '0.008164' },
    @{ Cell = 'E44'; Value = '-0.42%' },
    @{ Cell = 'D45'; Value = '0.3142' },
    @{ Cell = 'E45'; Value = '1.37%' },
    @{ Cell = 'D46'; Value = '0.00006653' },
    @{ Cell = 'E46'; Value = '-1.01%' },
    @{ Cell = 'D47'; Value = '0.00000000753' },
    @{ Cell = 'E47'; Value = '0.52%' },
    @{ Cell = 'D48'; Value = '0.2918' },
    @{ Cell = 'E48'; Value = '71.78%' },
    @{ Cell = 'D49'; Value = '0.004218' },
    @{ Cell = 'E49'; Value = '3.03%' },
    @{ Cell = 'D50'; Value = '0.00002109' },
    @{ Cell = 'E50'; Value = '0.52%' },
    @{ Cell = 'D51'; Value = '0.0002009' },
    @{ Cell = 'E51'; Value = '0.52%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
